$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (weeks): B1/D1 -> 16, C1/E1 -> 20
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 removed (Lichtwark deleted values); C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -3.9100618474287927
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): updated values
$ws.Range("B3").Value = -3.6976792365190221
$ws.Range("C3").Value = -1.112244729665953
$ws.Range("D3").Value = -6.8681887904253571
$ws.Range("E3").Value = 10.343436047236189

# Update active selection to match the new focused range
$ws.Range("B1:E3").Select()
